# Update analysis values for 103124 and 110724 image analysis on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: 1K_PFF_HTRA1_Hoechst_20x_03
$ws.Range("B4").Value = 43
$ws.Range("E4").Value = 6281
$ws.Range("F4").Value = 288

# Row 19: 1K_PFF_no_HTRA1_Hoechst_20x_08
$ws.Range("B19").Value = 11
$ws.Range("E19").Value = 1484

# Row 31: wt_PFF_HTRA1_Hoechst_20x_10
$ws.Range("B31").Value = 11
$ws.Range("E31").Value = 685
$ws.Range("F31").Value = 10
